$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("BonusPower", 0.03, 0,   22, 100, 0,   "lose"),
    @("BonusPower", 0.03, 0,   23, 100, 0,   "lose"),
    @("BonusPower", 2,    980, 98, 100, 200, "win"),
    @("BonusPower", 2,    560, 96, 100, 200, "win"),
    @("SkipBoss",   2,    960, 96, 100, 200, "win")
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
